$d = $word.ActiveDocument

# 1. Split the first paragraph in two: keep "...detector_energy_window_size "
#    in the first paragraph, and move what was the trailing (empty) run into a
#    brand-new paragraph that now reads "6. nonadiabatic_coupling".
$d.Content.Find.Execute(
    "detector_energy_window_size ", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "detector_energy_window_size ^p6. nonadiabatic_coupling", 2)

# Re-apply the (already-present) black color to just the trailing space so it
# keeps living in its own run, matching the original two-run split of
# paragraph 1 ("..size" run + " " run with explicit color formatting).
$spaceRange = $d.Range(115, 116)
$spaceRange.Font.Color = 0

# 2. maxdis line: add ", Franck_Condon_factor_cutoff" after "cutoff"
$d.Content.Find.Execute(
    "maxdis (for constructing anharmonic coupling)   cutoff  ", $true, $false,
    $false, $false, $false, $true, 1, $false,
    "maxdis (for constructing anharmonic coupling)   cutoff ,  Franck_Condon_factor_cutoff ", 2)

# 3. nmodes[i]  proptime[i] line: add two trailing spaces
$d.Content.Find.Execute(
    "nmodes[i]  proptime[i]", $true, $false, $false, $false, $false, $true, 1,
    $false, "nmodes[i]  proptime[i]  ", 2)

# 4. mfrequency  nmax line: append "  EV_coupling_V" (drop the old trailing space)
$d.Content.Find.Execute(
    "mfrequency  nmax ", $true, $false, $false, $false, $false, $true, 1,
    $false, "mfrequency  nmax  EV_coupling_V", 2)
